$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 17: quantity 12 -> 11 (C17 total quantity reduced) ---
# N17 = L17*K17 recalculates automatically via the shared formula.
$ws.Range("K17").Value = 11

# --- Append new BOM row 25: an additional 1uF 0603 capacitor line ---
$ws.Range("A25").Value = 25
$ws.Range("C25").Value = $ws.Range("C17").Value()
$ws.Range("D25").Value = "1uF"
$ws.Range("J25").Value = $ws.Range("J17").Value()
$ws.Range("K25").Value = 1

# The B/F/O cells on the new row stay empty, but still carry the same
# formatting as the rest of the table, so copy the formats over and then
# clear the (unused) content those source cells happened to carry.
$ws.Range("B17").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").ClearContents()
$ws.Range("B25").Font.Name = $ws.Range("B17").Font.Name()

$ws.Range("F17").Copy()
$ws.Range("F25").PasteSpecial(-4122)
$ws.Range("F25").ClearContents()

$ws.Range("O17").Copy()
$ws.Range("O25").PasteSpecial(-4122)
$ws.Range("O25").ClearContents()

# --- Update the view: scroll so column D / row 7 is the top-left cell and
#     select K26, mirroring where the user ended up after adding the row. ---
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 4
$ws.Range("K26").Select()

$wb.Save()
